# BOI11.xlsx update:
#  - The "Status" column (H) is recoded from a sales-pipeline scale
#    (presented/pending/declined/won) to a satisfaction scale
#    (Very bad/Bad/Good/Very good).
#  - The "Name" (B), "Rep" (C), "Manager" (D), "Product" (E) and
#    "Quantity" (F) columns are refreshed for rows 2-17.
#  - The final data row (18, account 729833 / Koepp Ltd / Monitor) is removed.
#  - The view is reset: selection moves to F3 and the old scrolled
#    top-left cell (A4) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status values (column H) -------------------------------------------
# Set these first, and in this exact order, so that the brand-new shared
# strings ("Bad", "Very bad", "Very good", "Good") are registered in the
# same order they appear in the target workbook.
$ws.Range("H4").Value  = "Bad"
$ws.Range("H5").Value  = "Very bad"
$ws.Range("H6").Value  = "Very good"
$ws.Range("H2").Value  = "Good"
$ws.Range("H3").Value  = "Good"
$ws.Range("H7").Value  = "Bad"
$ws.Range("H8").Value  = "Good"
$ws.Range("H9").Value  = "Bad"
$ws.Range("H10").Value = "Very bad"
$ws.Range("H11").Value = "Very good"
$ws.Range("H12").Value = "Good"
$ws.Range("H13").Value = "Bad"
$ws.Range("H14").Value = "Good"
$ws.Range("H15").Value = "Very good"
$ws.Range("H16").Value = "Very good"
$ws.Range("H17").Value = "Very bad"

# --- Remaining columns for each data row (2-17) --------------------------
$rows = @(
    @{ Row = 2;  Name = "Week Password";                   Rep = "Craig Booker";  Manager = "Juan Para";     Product = "CPU";         Quantity = 3 },
    @{ Row = 3;  Name = "CBTS overdue";                     Rep = "Craig Booker";  Manager = "Juan Para";     Product = "Software";    Quantity = 3 },
    @{ Row = 4;  Name = "PhishMe clicked";                  Rep = "Craig Booker";  Manager = "Juan Para";     Product = "Maintenance"; Quantity = 2 },
    @{ Row = 5;  Name = "Security Incident Involvement";    Rep = "Craig Booker";  Manager = "Juan Para";     Product = "CPU";         Quantity = 1 },
    @{ Row = 6;  Name = "Week Password";                    Rep = "Daniel Hilton"; Manager = "Felipe Fiorin"; Product = "CPU";         Quantity = 5 },
    @{ Row = 7;  Name = "CBTS overdue";                     Rep = "Daniel Hilton"; Manager = "Felipe Fiorin"; Product = "CPU";         Quantity = 2 },
    @{ Row = 8;  Name = "PhishMe clicked";                  Rep = "Daniel Hilton"; Manager = "Felipe Fiorin"; Product = "Software";    Quantity = 4 },
    @{ Row = 9;  Name = "Security Incident Involvement";    Rep = "John Smith";    Manager = "Felipe Fiorin"; Product = "Maintenance"; Quantity = 2 },
    @{ Row = 10; Name = "Week Password";                    Rep = "John Smith";    Manager = "Kevin Whelan";  Product = "CPU";         Quantity = 1 },
    @{ Row = 11; Name = "CBTS overdue";                     Rep = "Cedric Moss";   Manager = "Kevin Whelan";  Product = "CPU";         Quantity = 5 },
    @{ Row = 12; Name = "PhishMe clicked";                  Rep = "Cedric Moss";   Manager = "Kevin Whelan";  Product = "CPU";         Quantity = 4 },
    @{ Row = 13; Name = "Security Incident Involvement";    Rep = "Cedric Moss";   Manager = "Kevin Whelan";  Product = "Maintenance"; Quantity = 1 },
    @{ Row = 14; Name = "Week Password";                    Rep = "Cedric Moss";   Manager = "Maeve Morris";  Product = "Software";    Quantity = 4 },
    @{ Row = 15; Name = "CBTS overdue";                     Rep = "Wendy Yule";    Manager = "Maeve Morris";  Product = "Maintenance"; Quantity = 5 },
    @{ Row = 16; Name = "PhishMe clicked";                  Rep = "Wendy Yule";    Manager = "Maeve Morris";  Product = "CPU";         Quantity = 5 },
    @{ Row = 17; Name = "Security Incident Involvement";    Rep = "Wendy Yule";    Manager = "Maeve Morris";  Product = "CPU";         Quantity = 1 }
)

foreach ($r in $rows) {
    $ws.Range("B" + $r.Row).Value = $r.Name
    $ws.Range("C" + $r.Row).Value = $r.Rep
    $ws.Range("D" + $r.Row).Value = $r.Manager
    $ws.Range("E" + $r.Row).Value = $r.Product
    $ws.Range("F" + $r.Row).Value = $r.Quantity
}

# --- Drop the obsolete trailing row (18) ---------------------------------
$ws.Rows.Item(18).Delete()

# --- Reset the view: clear the scrolled top-left cell, select F3 --------
$ws.Range("F3").Select()
